# Rename the "wt" and "dcin5" worksheets to include "_log2_expression",
# matching the commit message: "Updated the input files in sixteen_tests
# to have strain_log2_expression instead of just strain".

$wb = $excel.ActiveWorkbook

$wsWt = $wb.Worksheets.Item("wt")
$wsWt.Name = "wt_log2_expression"

$wsDcin5 = $wb.Worksheets.Item("dcin5")
$wsDcin5.Name = "dcin5_log2_expression"

# Move the selection on the renamed "dcin5_log2_expression" sheet and make
# it the active (selected) tab.
$wsDcin5.Select()
$wsDcin5.Range("E41").Select()

# The previously active sheet ("optimization_parameters") loses the
# tabSelected flag automatically once another sheet is selected/activated.

# Scroll the tab strip so "wt_log2_expression" is the first visible tab.
$excel.ActiveWindow.ScrollWorkbookTabs(0, 2)
